# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Fri Sep 13 12:33:56 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.721.35'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.348.99'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '545.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.06'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.565'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +5.33%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.104'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.53'
$ws.Range('D10').Style = "Normal"
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.351'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '23.76'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = '2.768.11'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '57.656.93'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').Value = '2.350.48'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.98'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.67%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.28'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '328.12'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.87'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.63%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '63.26'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.166'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.18'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.31'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.61%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.75'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '171.08'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('D30').Value = '0.0₃0736'
$ws.Range('E30').Value = '  +0.77%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.11'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.31'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.53%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.13'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.22'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.58'
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.413'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +9.17%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '142.09'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.68%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '288.57'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.63'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0946'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0510'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.16%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.564'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '18.60'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('B47').Value = 'Polygon'
$ws.Range('C47').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.391'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.85%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0220'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.08'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.945'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.07%  '
